$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the "Project ID" column (A) with IDs parsed out of the
# "Project Name" column (B). Most are plain numbers; some projects use an
# alphanumeric code (S.., E.., P47) which must be written as text. Row 10
# (Project ID 9) is left untouched - it already had the correct value.

$ws.Range("A2").Value = 16
$ws.Range("A3").Value = "S11"
$ws.Range("A4").Value = "S16"
$ws.Range("A5").Value = 12
$ws.Range("A6").Value = 54
$ws.Range("A7").Value = "S17"
$ws.Range("A8").Value = "S02"
$ws.Range("A9").Value = 30
$ws.Range("A11").Value = 34
$ws.Range("A12").Value = "S07"
$ws.Range("A13").Value = 69
$ws.Range("A14").Value = 55
$ws.Range("A15").Value = "S03"
$ws.Range("A16").Value = "S15"
$ws.Range("A17").Value = "S08"
$ws.Range("A18").Value = "S04"
$ws.Range("A19").Value = "E03"
$ws.Range("A20").Value = "S13"
$ws.Range("A21").Value = 35
$ws.Range("A22").Value = 24
$ws.Range("A23").Value = 11
$ws.Range("A24").Value = 59

# Row 26 was typed before row 25 in the original edit (so "E01" grabbed a
# lower shared-string index than "01" even though it's the later row) -
# replicate that ordering here to match the shared string table exactly.
$ws.Range("A26").Value = "E01"

# "01" must keep its leading zero, so force a text number format before
# assigning the string value.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "01"

$ws.Range("A27").Value = 7
$ws.Range("A28").Value = 32
$ws.Range("A29").Value = "S09"
$ws.Range("A30").Value = 43
$ws.Range("A31").Value = 17
$ws.Range("A32").Value = 45
$ws.Range("A33").Value = 50
$ws.Range("A34").Value = 10
$ws.Range("A35").Value = 13
$ws.Range("A36").Value = 36
$ws.Range("A37").Value = "S01"
$ws.Range("A38").Value = 46
$ws.Range("A39").Value = "S18"
$ws.Range("A40").Value = 44
$ws.Range("A41").Value = 21
$ws.Range("A42").Value = "P47"
$ws.Range("A43").Value = "E02"
$ws.Range("A44").Value = 15
$ws.Range("A45").Value = "E05"

# Move the active selection down to the bottom of the now-filled column.
$ws.Range("A45").Select()
